# Update the "Template.xlsx" Summary sheet to match the updated
# safedata_validator documentation template.
#
# Net effect (per the target OOXML):
#  - shared string "SAFE Project ID" is retired and replaced by "Project ID"
#  - four new label rows are introduced: "Access conditions", "Permit type",
#    "Permit authority" and "Permit number"
#  - the Summary sheet labels are reordered/extended from 27 to 31 rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

$labels = @(
    "Project ID",
    "Title",
    "Description",
    "Access status",
    "Embargo date",
    "Access conditions",
    "Author name",
    "Author email",
    "Author affiliation",
    "Author ORCID",
    "Worksheet name",
    "Worksheet title",
    "Worksheet description",
    "Worksheet external file",
    "Keywords",
    "External file",
    "External file description",
    "Publication DOI",
    "Funding body",
    "Funding type",
    "Funding reference",
    "Funding link",
    "Permit type",
    "Permit authority",
    "Permit number",
    "Start date",
    "End date",
    "North",
    "South",
    "East",
    "West"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $labels[$i]
}

# Keep the active sheet selection in line with the authored workbook state.
$ws.Range("C21:C22").Select() | Out-Null
